$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed cryptocurrency price/volume data (cryptos list update).
# Row 50/51 also swap Coin/Link (ThetaToken <-> BabyDogeCoin) per the source feed re-ranking.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.600.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.88%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.005.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.37%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5021"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4261"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.56%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.79"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09123"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.87%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.128"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.83%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.126"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.550"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.966.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.007"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001120"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06650"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.62%  "

$ws.Range("E21").Value = "  +0.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.987"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.69%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.582.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.69%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.268"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.559"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.341"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.49%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.057"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.92%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.588"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -10.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09956"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.874"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.71%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.781"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.664"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02476"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.90%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06392"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.306"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6574"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.74%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2075"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.91%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.006"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6362"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.97%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.219"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.279"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.534"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07012"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.31%  "

$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.140"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.53%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000322"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.60%  "
